$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on A/B columns so numeric-looking strings stay as text
$ws.Range("A5").NumberFormat = "@"
$ws.Range("B2:B5").NumberFormat = "@"

# New row 5 - A value first, to mirror shared-string insertion order
$ws.Range("A5").Value = "4"

# Update existing rows: B column text values and C column numeric values change.
$ws.Range("B2").Value = "2016"
$ws.Range("C2").Value = 1484258.0

$ws.Range("B3").Value = "2018"
$ws.Range("C3").Value = 1.0

$ws.Range("B4").Value = "2106"
$ws.Range("C4").Value = 6.0

$ws.Range("B5").Value = "2916"
$ws.Range("C5").Value = 1.0
